$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per the diff (NATMI TPM recompute: Wnt5a-Fzd1)
$ws.Range("G2").Value = 0.02354566666666667
$ws.Range("H2").Value = 0.07063700000000001
$ws.Range("I2").Value = 0.002815555392485919
$ws.Range("J2").Value = 0.002815555392485918
$ws.Range("M2").Value = 0.8596446666666667
$ws.Range("N2").Value = 2.578934
$ws.Range("O2").Value = 0.05286426382906832
$ws.Range("P2").Value = 0.05286426382906832
$ws.Range("Q2").Value = 0.02024090677311112
$ws.Range("R2").Value = 0.182168160958
$ws.Range("S2").Value = 0.0001488422630937316
$ws.Range("T2").Value = 0.0001488422630937316
$ws.Range("G3").Value = 0.02354566666666667
$ws.Range("H3").Value = 0.07063700000000001
$ws.Range("I3").Value = 0.002815555392485919
$ws.Range("J3").Value = 0.002815555392485918
$ws.Range("O3").Value = 0.6417658132713033
$ws.Range("P3").Value = 0.6417658132713032
$ws.Range("Q3").Value = 0.245722177057
$ws.Range("R3").Value = 2.211499593513
$ws.Range("S3").Value = 0.001806927196269129
$ws.Range("T3").Value = 0.001806927196269129
$ws.Range("G4").Value = 0.02354566666666667
$ws.Range("H4").Value = 0.07063700000000001
$ws.Range("I4").Value = 0.002815555392485919
$ws.Range("J4").Value = 0.002815555392485918
$ws.Range("O4").Value = 0.3053699228996285
$ws.Range("P4").Value = 0.3053699228996284
$ws.Range("Q4").Value = 0.1169214076395556
$ws.Range("R4").Value = 1.052292668756
$ws.Range("S4").Value = 0.0008597859331230583
$ws.Range("T4").Value = 0.0008597859331230579
$ws.Range("I5").Value = 0.9868456480383168
$ws.Range("J5").Value = 0.9868456480383166
$ws.Range("M5").Value = 0.8596446666666667
$ws.Range("N5").Value = 2.578934
$ws.Range("O5").Value = 0.05286426382906832
$ws.Range("P5").Value = 0.05286426382906832
$ws.Range("Q5").Value = 7.094390973341111
$ws.Range("R5").Value = 63.84951876007
$ws.Range("S5").Value = 0.05216886869646547
$ws.Range("T5").Value = 0.05216886869646546
$ws.Range("I6").Value = 0.9868456480383168
$ws.Range("J6").Value = 0.9868456480383166
$ws.Range("O6").Value = 0.6417658132713033
$ws.Range("P6").Value = 0.6417658132713032
$ws.Range("Q6").Value = 86.12505429740499
$ws.Range("S6").Value = 0.6333237998865567
$ws.Range("T6").Value = 0.6333237998865565
$ws.Range("I7").Value = 0.9868456480383168
$ws.Range("J7").Value = 0.9868456480383166
$ws.Range("O7").Value = 0.3053699228996285
$ws.Range("P7").Value = 0.3053699228996284
$ws.Range("S7").Value = 0.3013529794552947
$ws.Range("T7").Value = 0.3013529794552946
$ws.Range("I8").Value = 0.0103387965691973
$ws.Range("J8").Value = 0.0103387965691973
$ws.Range("M8").Value = 0.8596446666666667
$ws.Range("N8").Value = 2.578934
$ws.Range("O8").Value = 0.05286426382906832
$ws.Range("P8").Value = 0.05286426382906832
$ws.Range("Q8").Value = 0.07432516442822222
$ws.Range("R8").Value = 0.668926479854
$ws.Range("S8").Value = 0.0005465528695091126
$ws.Range("T8").Value = 0.0005465528695091126
$ws.Range("I9").Value = 0.0103387965691973
$ws.Range("J9").Value = 0.0103387965691973
$ws.Range("O9").Value = 0.6417658132713033
$ws.Range("P9").Value = 0.6417658132713032
$ws.Range("Q9").Value = 0.9022985688409997
$ws.Range("S9").Value = 0.006635086188477468
$ws.Range("T9").Value = 0.006635086188477466
$ws.Range("I10").Value = 0.0103387965691973
$ws.Range("J10").Value = 0.0103387965691973
$ws.Range("O10").Value = 0.3053699228996285
$ws.Range("P10").Value = 0.3053699228996284
$ws.Range("S10").Value = 0.003157157511210724
$ws.Range("T10").Value = 0.003157157511210723
